$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 340, shifting existing rows 340-404 down to 341-405
$ws.Rows.Item(340).Insert()

# Populate the new row 340 with the new weekly data point
$ws.Cells.Item(340, 1).Value2 = 7
$ws.Cells.Item(340, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(340, 3).Value2 = "Ñuble"
$ws.Cells.Item(340, 4).Value2 = 45015
$ws.Cells.Item(340, 5).Value2 = 16
$ws.Cells.Item(340, 6).Value2 = 100114013
$ws.Cells.Item(340, 7).Value2 = "Zanahoria"
$ws.Cells.Item(340, 8).Value2 = "Sin especificar"
$ws.Cells.Item(340, 9).Value2 = "Primera"
$ws.Cells.Item(340, 10).Value2 = 80
$ws.Cells.Item(340, 11).Value2 = 7000
$ws.Cells.Item(340, 12).Value2 = 7000
$ws.Cells.Item(340, 13).Value2 = 7000
$ws.Cells.Item(340, 14).Value2 = "$/saco 20 kilos"
$ws.Cells.Item(340, 15).Value2 = "Región de Ñuble"
$ws.Cells.Item(340, 16).Value2 = 350
$ws.Cells.Item(340, 17).Value2 = 20
$ws.Cells.Item(340, 18).Value2 = "Hortaliza"

"Done"
